$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (row 11, column B)
$ws.Range("B11").Value = 5

# Update total correct count (row 12, column B)
$ws.Range("B12").Value = 95

# Update correct/total marks text (row 12, column E)
$ws.Range("E12").Value = "95/140"
